$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '317.07'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '-3.19%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '41.90'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '-5.55%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.202'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '1.64%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.08079'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '-3.53%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '4.373'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '-1.59%'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.752'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '-9.54%'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.9301'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '-4.54%'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '-0.23%'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1859'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '-2.27%'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.09321'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '-3.67%'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.04569'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-1.17%'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.387'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-16.88%'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.1054'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-0.56%'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.001294'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '0.31%'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.005913'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '1.92%'
$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.356'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '-1.41%'
$ws.Range("B18").Value = 'BTSEToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.547'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '1.04%'
$ws.Range("B19").Value = 'BitpandaEcosystemToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.3397'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '1.08%'
$ws.Range("B20").Value = 'ProBitToken'
$ws.Range("C20").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.1383'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '1.47%'
$ws.Range("B21").Value = 'ZBToken'
$ws.Range("C21").Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.2550'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-1.24%'
$ws.Range("B22").Value = 'CoinExToken'
$ws.Range("C22").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.04166'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '0.00%'
$ws.Range("B23").Value = 'BitKan'
$ws.Range("C23").Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.001244'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '0.39%'
$ws.Range("B24").Value = 'HotbitToken'
$ws.Range("C24").Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.004324'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '-2.26%'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0001224'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '-6.28%'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0002987'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '0.05%'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02577'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '-5.34%'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05448'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '-3.09%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.008043'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '2.46%'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '-1.31%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.007571'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '2.75%'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002086'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-1.44%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.008272'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '4.37%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.3143'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '-10.25%'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00006780'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '-1.89%'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000752'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '0.05%'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.003392'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '-3.30%'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.004111'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '16.20%'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00002106'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '0.05%'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0002005'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '0.05%'
